$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 37241.9
$ws.Range("I62").Value = 999
$ws.Range("J62").Value = 41268.89
$ws.Range("K62").Value = 999
$ws.Range("L62").Value = 41268.89
$ws.Range("N62").Value = -42516.89
$ws.Range("M62").Value = -375

# Row 65
$ws.Range("H65").Value = 37241.9
$ws.Range("I65").Value = 999
$ws.Range("J65").Value = 41268.89
$ws.Range("K65").Value = 4995
$ws.Range("L65").Value = 206344.45
$ws.Range("N65").Value = -212584.45
$ws.Range("M65").Value = -1875

# Row 127
$ws.Range("H127").Value = 2753.24
$ws.Range("J127").Value = 3787.8572
$ws.Range("L127").Value = 11363.5716
$ws.Range("N127").Value = -21283.5716

# Row 129
$ws.Range("H129").Value = 628.0909
$ws.Range("I129").Value = 628.0909
$ws.Range("K129").Value = 1884.2727
$ws.Range("M129").Value = 3115.7273

# Row 132
$ws.Range("H132").Value = 1831.4688
$ws.Range("I132").Value = 1641.6207
$ws.Range("K132").Value = 4924.8621
$ws.Range("M132").Value = -2394.8621

# Row 138
$ws.Range("H138").Value = 5558.9395
$ws.Range("I138").Value = 1965.6666
$ws.Range("J138").Value = 6126.2983
$ws.Range("K138").Value = 5896.9998
$ws.Range("L138").Value = 18378.8949
$ws.Range("M138").Value = -756.9997999999996
$ws.Range("N138").Value = -28658.8949

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1609267.4
$ws.Range("I32").Value = 1694700.6
$ws.Range("K32").Value = 1694700.6
$ws.Range("M32").Value = -1694413.6

# Row 45
$ws.Range("H45").Value = 4779.875
$ws.Range("I45").Value = 3931.7368
$ws.Range("K45").Value = 3931.7368
$ws.Range("M45").Value = -3554.7368

# Row 46
$ws.Range("H46").Value = 5787.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5787.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5787.75
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -6425.75

# Row 61
$ws.Range("H61").Value = 8645.789000000001
$ws.Range("I61").Value = 783.5
$ws.Range("K61").Value = 783.5
$ws.Range("M61").Value = -571.5

# Row 74
$ws.Range("H74").Value = 37940.895
$ws.Range("I74").Value = 51436.55
$ws.Range("J74").Value = 4201.75
$ws.Range("K74").Value = 51436.55
$ws.Range("L74").Value = 4201.75
$ws.Range("M74").Value = -50562.55
$ws.Range("N74").Value = -5949.75

# Row 77
$ws.Range("H77").Value = 37940.895
$ws.Range("I77").Value = 51436.55
$ws.Range("J77").Value = 4201.75
$ws.Range("K77").Value = 257182.75
$ws.Range("L77").Value = 21008.75
$ws.Range("M77").Value = -252814.75
$ws.Range("N77").Value = -29744.75

# Row 132
$ws.Range("H132").Value = 4757.672
$ws.Range("I132").Value = 3531.641
$ws.Range("K132").Value = 10594.923
$ws.Range("M132").Value = -8064.923000000001

# Row 136
$ws.Range("H136").Value = 8645.789000000001
$ws.Range("I136").Value = 783.5
$ws.Range("K136").Value = 2350.5
$ws.Range("M136").Value = 199.5

$ws = $wb.Worksheets.Item("BSM")
# Row 68
$ws.Range("H68").Value = 43000
$ws.Range("J68").Value = 43000
$ws.Range("L68").Value = 43000
$ws.Range("N68").Value = -44622

# Row 71
$ws.Range("H71").Value = 43000
$ws.Range("J71").Value = 43000
$ws.Range("L71").Value = 129000
$ws.Range("N71").Value = -137112

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 255.71428
$ws.Range("I22").Value = 197.25
$ws.Range("J22").Value = 333.66666
$ws.Range("K22").Value = 197.25
$ws.Range("L22").Value = 333.66666
$ws.Range("M22").Value = 152.75
$ws.Range("N22").Value = -1033.66666

# Row 31
$ws.Range("H31").Value = 5417.4175
$ws.Range("I31").Value = 2543.923
$ws.Range("K31").Value = 2543.923
$ws.Range("M31").Value = -2248.923

# Row 34
$ws.Range("H34").Value = 5417.4175
$ws.Range("I34").Value = 2543.923
$ws.Range("K34").Value = 2543.923
$ws.Range("M34").Value = -2341.923

# Row 58
$ws.Range("H58").Value = 16674987
$ws.Range("I58").Value = 41669644
$ws.Range("J58").Value = 11882.389
$ws.Range("K58").Value = 41669644
$ws.Range("L58").Value = 11882.389
$ws.Range("M58").Value = -41669441
$ws.Range("N58").Value = -12288.389

# Row 116
$ws.Range("H116").Value = 77244
$ws.Range("J116").Value = 77244
$ws.Range("L116").Value = 77244
$ws.Range("N116").Value = -86422

# Row 119
$ws.Range("H119").Value = 45000
$ws.Range("J119").Value = 45000
$ws.Range("L119").Value = 45000
$ws.Range("N119").Value = -54676

# Row 132
$ws.Range("H132").Value = 5698.089
$ws.Range("I132").Value = 3011.4814
$ws.Range("J132").Value = 9728
$ws.Range("K132").Value = 9034.4442
$ws.Range("L132").Value = 29184
$ws.Range("M132").Value = -6504.4442
$ws.Range("N132").Value = -34244

# Row 134
$ws.Range("H134").Value = 3998.9734
$ws.Range("I134").Value = 1787.7291
$ws.Range("J134").Value = 7930.074
$ws.Range("K134").Value = 5363.1873
$ws.Range("L134").Value = 23790.222
$ws.Range("M134").Value = -2828.1873
$ws.Range("N134").Value = -28860.222

# Row 136
$ws.Range("H136").Value = 16674987
$ws.Range("I136").Value = 41669644
$ws.Range("J136").Value = 11882.389
$ws.Range("K136").Value = 125008932
$ws.Range("L136").Value = 35647.167
$ws.Range("M136").Value = -125006382
$ws.Range("N136").Value = -40747.167

# Row 141
$ws.Range("H141").Value = 76314.45
$ws.Range("J141").Value = 76314.45
$ws.Range("L141").Value = 76314.45
$ws.Range("N141").Value = -86674.45

$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

# Row 109
$ws.Range("H109").Value = 154520880
$ws.Range("I109").Value = 168990800
$ws.Range("K109").Value = 506972400
$ws.Range("M109").Value = -506971360

# Row 125
$ws.Range("H125").Value = 125004500
$ws.Range("I125").Value = 250002000
$ws.Range("K125").Value = 750006000
$ws.Range("M125").Value = -750001080

# Row 138
$ws.Range("H138").Value = 50760.773
$ws.Range("I138").Value = 60042.555
$ws.Range("K138").Value = 180127.665
$ws.Range("M138").Value = -174987.665

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 6499.8335

# Row 122
$ws.Range("H122").Value = 1964655.6
$ws.Range("I122").Value = 3027515.2
$ws.Range("J122").Value = 2453.1538
$ws.Range("K122").Value = 9082545.600000001
$ws.Range("L122").Value = 7359.4614
$ws.Range("M122").Value = -9080095.600000001
$ws.Range("N122").Value = -12259.4614

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5012.048
$ws.Range("I40").Value = 3283.5334
$ws.Range("K40").Value = 3283.5334
$ws.Range("M40").Value = -3147.5334

# Row 46
$ws.Range("H46").Value = 5295034.5
$ws.Range("J46").Value = 5852143.5
$ws.Range("L46").Value = 5852143.5
$ws.Range("N46").Value = -5852519.5

# Row 76
$ws.Range("H76").Value = 14816.667
$ws.Range("J76").Value = 14816.667
$ws.Range("L76").Value = 14816.667
$ws.Range("N76").Value = -15492.667

# Row 79
$ws.Range("H79").Value = 14816.667
$ws.Range("J79").Value = 14816.667
$ws.Range("L79").Value = 14816.667
$ws.Range("N79").Value = -17156.667

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 2719.4
$ws.Range("J41").Value = 2719.4
$ws.Range("L41").Value = 2719.4
$ws.Range("N41").Value = -3499.4

# Row 136
$ws.Range("H136").Value = 31286468
$ws.Range("I136").Value = 71429790
$ws.Range("K136").Value = 214289370
$ws.Range("M136").Value = -214286820
